$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text format so numeric-looking strings (dates-like prices, percentages)
# are preserved exactly as text, matching the original inline-string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.031.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6277"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07579"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.84"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07648"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.830.36"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.955"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.46"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009519"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +10.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.990"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.015.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "225.10"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.207"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.424"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1364"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.063"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.034"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05210"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7321"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.591"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.269.53"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.757"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01788"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.533"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8907"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.51"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.977.60"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.77"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5112"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3979"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.861"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05751"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.635"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.35%  "
